$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the price-list date in A1 (serial date: 2024-05-03 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update prices in column D for rows 29-36
$ws.Range("D29").Value = 2137
$ws.Range("D30").Value = 2322
$ws.Range("D31").Value = 1638
$ws.Range("D32").Value = 1726
$ws.Range("D33").Value = 1750
$ws.Range("D34").Value = 1861
$ws.Range("D35").Value = 1876
$ws.Range("D36").Value = 2052
